# Manually add ICDC TCL01 Test for Prasanna
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# ---- long query / filename strings used below (here-strings preserve content verbatim) ----
$s9 = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
   count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$s10 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
MATCH (samp:sample)-->(c)
WHERE  samp.specific_sample_pathology in ['Lymphoma']
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$s13 = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$s14 = @'
TC03_Canine_Filter_SamplePatho-Lymohoma_Neo4jData.xlsx
'@

$s15 = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
MATCH (samp:sample)-->(c)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$s16 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE  samp.specific_sample_pathology IN ['Lymphoma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

$s17 = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (f)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Lymphoma']
 MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 100
'@

# ---- Row 2 (CasesTab) ----
$ws.Range("B2").Value2 = $s10
$ws.Range("C2").Value2 = $s9
$ws.Range("D2").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_Neo4jData.xlsx"
$ws.Range("E2").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_WebData.xlsx"

# ---- Row 3 (SamplesTab) ----
$ws.Range("B3").Value2 = $s16
$ws.Range("C3").Value2 = $s13
$ws.Range("D3").Value2 = $s14
$ws.Range("E3").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_WebData.xlsx"

# ---- Row 4 (FilesTab) ----
$ws.Range("A4").Value2 = "FilesTab"
$ws.Range("B4").Value2 = $s17
$ws.Range("C4").Value2 = $s13
$ws.Range("D4").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_Neo4jData.xlsx"
$ws.Range("E4").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_WebData.xlsx"

# ---- Row 5 (new StudyFilesTab) ----
$ws.Range("A5").Value2 = "StudyFilesTab"
$ws.Range("B5").Value2 = $s15
$ws.Range("C5").Value2 = $s13
$ws.Range("D5").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_Neo4jData.xlsx"
$ws.Range("E5").Value2 = "TC03_Canine_Filter_SamplePatho-Lymphoma_WebData.xlsx"

# ---- Formatting ----
# Existing wrap style (s=1) already applied to B2/C2 originally; keep B2 as-is (wrap only)
$ws.Range("B2").WrapText = $true

# C2, C3, C4 get wrap + vertical-center (new style s=2)
$ws.Range("C2:C4").WrapText = $true
$ws.Range("C2:C4").VerticalAlignment = -4108  # xlCenter

# B3, B4, B5 get a bigger font (18pt) + wrap (new style s=3)
$ws.Range("B3:B5").WrapText = $true
$ws.Range("B3:B5").Font.Size = 18

# C5 also uses the bigger-font style (s=3)
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Font.Size = 18

# ---- Row heights (content grew substantially; Excel auto row-height does not
#      recompute in this headless host, so set explicitly to match real Excel output) ----
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# ---- Selection / scroll position ----
$ws.Range("C5").Select()

